$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  "D2" = "257.47"
  "E2" = "1.14%"
  "D3" = "26.99"
  "E3" = "-3.93%"
  "D4" = "4.743"
  "E4" = "-10.93%"
  "D5" = "0.05968"
  "E5" = "2.01%"
  "D6" = "6.680"
  "E6" = "-0.50%"
  "D7" = "0.8705"
  "E7" = "0.45%"
  "D8" = "0.9465"
  "E8" = "3.78%"
  "D9" = "0.1406"
  "E9" = "-1.29%"
  "E10" = "4.74%"
  "D11" = "0.07191"
  "E11" = "0.26%"
  "D12" = "0.03173"
  "E12" = "-0.25%"
  "D13" = "0.09237"
  "E13" = "0.17%"
  "D14" = "0.001551"
  "E14" = "-0.18%"
  "D15" = "0.0006124"
  "E15" = "0.90%"
  "D16" = "0.005994"
  "E16" = "3.24%"
  "D17" = "3.499"
  "E17" = "-0.03%"
  "D18" = "3.178"
  "E18" = "-1.60%"
  "D20" = "0.3099"
  "E20" = "-2.21%"
  "E21" = "-0.65%"
  "D22" = "3.815"
  "E22" = "7.40%"
  "D23" = "0.04224"
  "E23" = "1.32%"
  "E24" = "0.17%"
  "E25" = "0.10%"
  "D26" = "0.004500"
  "D27" = "0.0001201"
  "E27" = "0.11%"
  "E28" = "-22.90%"
  "D40" = "0.03823"
  "E40" = "-0.43%"
  "D41" = "0.006205"
  "E41" = "62.13%"
  "D42" = "0.1100"
  "E42" = "-0.12%"
  "D43" = "0.002253"
  "E43" = "-5.68%"
  "D44" = "0.01060"
  "E44" = "-3.21%"
  "D45" = "0.00005504"
  "E45" = "4.98%"
  "E46" = "0.08%"
  "D47" = "0.1092"
  "E47" = "21.52%"
  "D48" = "0.002277"
  "E48" = "5.59%"
  "D49" = "0.00002101"
  "E49" = "0.08%"
  "D50" = "0.0002001"
  "E50" = "0.08%"
}

foreach ($addr in $updates.Keys) {
  $cell = $ws.Range($addr)
  $cell.NumberFormat = "@"
  $cell.Value = $updates[$addr]
  $cell.Style = "Normal"
}
